# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps for the handback run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# zh-cn: row 3 corresponds to the 9a3c42df... file that was just handed back.
$wsZhCn.Range("E3").Value = "2016-03-20 00:36:57"
$wsZhCn.Range("H3").Value = "2016-03-20 00:37:16"

# de-de: row 3 corresponds to the 9a3c42df... file that was just handed back.
$wsDeDe.Range("E3").Value = "2016-03-20 00:36:59"
$wsDeDe.Range("H3").Value = "2016-03-20 00:37:21"
